# "Best practises update + testrun data"
# Updates the Input sheet's multi-destination dates and refreshes the
# Output sheet with a new test run's captured results.

$wb = $excel.ActiveWorkbook
$input = $wb.Worksheets.Item("Input")
$output = $wb.Worksheets.Item("Output")

# --- Input sheet: multi destination dates (row 4) ---
$input.Range("M4").Value = 27022022
$input.Range("O4").Value = 28022022

# --- Output sheet: new test-run results ---

# Row 2
$output.Range("A2").Value = "28/01/2022 12:46:35 pm"
$output.Range("D2").Value = "₹1,21,520"
$output.Range("E2").Value = "₹8,364"
$output.Range("G2").Value = "₹1,29,894"

# Row 3
$output.Range("A3").Value = "28/01/2022 12:46:55 pm"
$output.Range("D3").Value = "₹26,395"
$output.Range("E3").Value = "₹3,390"
$output.Range("G3").Value = "₹29,795"

# Row 5
$output.Range("A5").Value = "28/01/2022 12:48:09 pm"

# Row 6
$output.Range("A6").Value = "28/01/2022 12:48:36 pm"
$output.Range("D6").Value = "₹1,44,310"
$output.Range("E6").Value = "₹10,344"
$output.Range("G6").Value = "₹1,54,664"

# Row 8
$output.Range("A8").Value = "28/01/2022 12:49:19 pm"

# Row 9
$output.Range("A9").Value = "28/01/2022 12:49:30 pm"

# Row 10
$output.Range("A10").Value = "28/01/2022 12:49:34 pm"

# Row 11
$output.Range("A11").Value = "28/01/2022 12:49:44 pm"

# --- Selection on Input sheet moves to O12 ---
$input.Range("O12").Select()
